# "Add q drop to master db"
#
# Adds a new "% of Q Drop's" column (column I) to the grade-distribution
# sheet, to the right of the existing "% of F's" column (H), and fills in
# the per-course/per-professor Q-drop percentages.
#
# NOTE: in this workbook the "% of ..." figures are stored as plain text
# (e.g. "3.85%"), not as numbers with a percentage number format, so each
# new cell is forced to Text format before its value is written -
# otherwise Excel would parse "0.00%" as the number 0 displayed with a
# "0.00%" format, which would not match the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "% of Q Drop's"

# Data rows (course/professor rows only; section-header rows like A2, A6, ...
# have no grade-distribution figures and therefore no Q-drop value either)
$qdrop = @{
    3  = "0.00%"
    4  = "1.69%"
    7  = "0.00%"
    10 = "0.00%"
    13 = "0.69%"
    16 = "0.00%"
    19 = "2.86%"
    22 = "0.00%"
    23 = "0.00%"
    26 = "0.00%"
    29 = "0.00%"
    30 = "0.00%"
    31 = "0.00%"
    34 = "0.00%"
    37 = "0.00%"
    40 = "0.00%"
    41 = "0.00%"
    42 = "0.00%"
    45 = "0.00%"
    48 = "0.00%"
}

foreach ($row in $qdrop.Keys) {
    $cell = $ws.Range("I$row")
    $cell.NumberFormat = "@"
    $cell.Value = $qdrop[$row]
}
